$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=78; A=4; B="Feria Lagunitas de Puerto Montt"; C="Los Lagos"; D=44656; E=10; F="Fruta"; G=100104; H="Frutos de pepita"; I=100104003; J="Membrillo"; K="Champion"; L="Especial";  M=150; N=19000; O=19000; P=19000; Q="$/caja 18 kilos granel"; Origen="Región de O'Higgins"; S=1056; T=18 },
    @{ Row=79; A=4; B="Feria Lagunitas de Puerto Montt"; C="Los Lagos"; D=44656; E=10; F="Fruta"; G=100104; H="Frutos de pepita"; I=100104003; J="Membrillo"; K="Champion"; L="Primera";   M=150; N=16000; O=16000; P=16000; Q="$/caja 18 kilos granel"; Origen="Región de O'Higgins"; S=889;  T=18 },
    @{ Row=80; A=4; B="Feria Lagunitas de Puerto Montt"; C="Los Lagos"; D=44656; E=10; F="Fruta"; G=100104; H="Frutos de pepita"; I=100104003; J="Membrillo"; K="Champion"; L="Segunda";   M=150; N=15000; O=15000; P=15000; Q="$/caja 18 kilos granel"; Origen="Región de O'Higgins"; S=833;  T=18 }
)

foreach ($row in $rows) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.Origen
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
}
